$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '45.289.90'
Set-TextValue $ws.Range("E2") '  +4.21%  '

Set-TextValue $ws.Range("D3") '2.424.74'
Set-TextValue $ws.Range("E3") '  +0.39%  '

Set-TextValue $ws.Range("E4") '  -0.14%  '

Set-TextValue $ws.Range("D5") '317.90'
Set-TextValue $ws.Range("E5") '  +3.75%  '

Set-TextValue $ws.Range("D6") '102.38'
Set-TextValue $ws.Range("E6") '  +5.47%  '

Set-TextValue $ws.Range("E7") '  +1.73%  '

Set-TextValue $ws.Range("D8") '0.999'
Set-TextValue $ws.Range("E8") '  -0.12%  '

Set-TextValue $ws.Range("D9") '0.529'
Set-TextValue $ws.Range("E9") '  +7.99%  '

Set-TextValue $ws.Range("D10") '35.59'
Set-TextValue $ws.Range("E10") '  +2.11%  '

Set-TextValue $ws.Range("E11") '  +0.88%  '

Set-TextValue $ws.Range("E12") '  -2.00%  '

Set-TextValue $ws.Range("D13") '18.10'
Set-TextValue $ws.Range("E13") '  -1.85%  '

Set-TextValue $ws.Range("E14") '  +2.24%  '

Set-TextValue $ws.Range("D15") '2.805.05'
Set-TextValue $ws.Range("E15") '  +0.63%  '

Set-TextValue $ws.Range("D16") '2.437.28'

Set-TextValue $ws.Range("D17") '0.842'
Set-TextValue $ws.Range("E17") '  +2.14%  '

Set-TextValue $ws.Range("D18") '45.187.62'
Set-TextValue $ws.Range("E18") '  +3.85%  '

Set-TextValue $ws.Range("D19") '12.22'
Set-TextValue $ws.Range("E19") '  +1.39%  '

Set-TextValue $ws.Range("E20") '  -0.90%  '

Set-TextValue $ws.Range("E21") '  +2.45%  '

Set-TextValue $ws.Range("D22") '68.75'

Set-TextValue $ws.Range("D23") '244.02'
Set-TextValue $ws.Range("E23") '  +2.70%  '

Set-TextValue $ws.Range("E24") '  +0.19%  '

Set-TextValue $ws.Range("D25") '2.50'
Set-TextValue $ws.Range("E25") '  +1.90%  '

Set-TextValue $ws.Range("E26") '  -0.04%  '

Set-TextValue $ws.Range("D27") '25.54'
Set-TextValue $ws.Range("E27") '  +2.50%  '

Set-TextValue $ws.Range("E28") '  -0.53%  '

Set-TextValue $ws.Range("E29") '  +1.64%  '

Set-TextValue $ws.Range("D30") '49.13'
Set-TextValue $ws.Range("E30") '  +2.18%  '

Set-TextValue $ws.Range("D31") '32.88'
Set-TextValue $ws.Range("E31") '  +2.14%  '

Set-TextValue $ws.Range("B32") 'Kaspa'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D32") '0.126'
Set-TextValue $ws.Range("E32") '  +5.28%  '

Set-TextValue $ws.Range("B33") 'Celestia'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D33") '20.24'
Set-TextValue $ws.Range("E33") '  +10.02%  '

Set-TextValue $ws.Range("D34") '5.20'
Set-TextValue $ws.Range("E34") '  +1.57%  '

Set-TextValue $ws.Range("E35") '  +0.18%  '

Set-TextValue $ws.Range("D36") '0.0763'
Set-TextValue $ws.Range("E36") '  +1.65%  '

Set-TextValue $ws.Range("E37") '  -0.89%  '

Set-TextValue $ws.Range("D38") '4.45'
Set-TextValue $ws.Range("E38") '  +1.79%  '

Set-TextValue $ws.Range("B39") 'Monero'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D39") '125.78'
Set-TextValue $ws.Range("E39") '  -5.17%  '

Set-TextValue $ws.Range("B40") 'LidoDAOToken'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D40") '2.85'
Set-TextValue $ws.Range("E40") '  -2.23%  '

Set-TextValue $ws.Range("B41") 'WEMIXToken'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D41") '2.22'
Set-TextValue $ws.Range("E41") '  -2.61%  '

Set-TextValue $ws.Range("B42") 'Stellar'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D42") '0.109'
Set-TextValue $ws.Range("E42") '  +1.07%  '

Set-TextValue $ws.Range("D43") '20.70'
Set-TextValue $ws.Range("E43") '  -0.94%  '

Set-TextValue $ws.Range("D44") '0.0290'
Set-TextValue $ws.Range("E44") '  +2.68%  '

Set-TextValue $ws.Range("D45") '1.933.78'
Set-TextValue $ws.Range("E45") '  -0.58%  '

Set-TextValue $ws.Range("E46") '  -2.99%  '

Set-TextValue $ws.Range("E47") '  +3.39%  '

Set-TextValue $ws.Range("E48") '  +16.44%  '

Set-TextValue $ws.Range("E49") '  -2.45%  '

Set-TextValue $ws.Range("D50") '76.38'
Set-TextValue $ws.Range("E50") '  +5.69%  '

Set-TextValue $ws.Range("D51") '53.85'
Set-TextValue $ws.Range("E51") '  +2.40%  '
